# Weekly update: insert this week's two new price rows (Primera / Segunda)
# at the top of the data block (row 15), pushing the existing history down.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows before row 15 (existing rows 15-53 shift to 17-55)
$ws.Range("A15:A16").EntireRow.Insert()

# New row 15 - "Primera" quality, this week's data
$ws.Cells.Item(15, 1).Value = 11
$ws.Cells.Item(15, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(15, 3).Value = "Bíobío"
$ws.Cells.Item(15, 4).Value = 44742
$ws.Cells.Item(15, 5).Value = 8
$ws.Cells.Item(15, 6).Value = 100112043
$ws.Cells.Item(15, 7).Value = "Pepino dulce"
$ws.Cells.Item(15, 8).Value = "Cultivar IV Región"
$ws.Cells.Item(15, 9).Value = "Primera"
$ws.Cells.Item(15, 10).Value = 100
$ws.Cells.Item(15, 11).Value = 14000
$ws.Cells.Item(15, 12).Value = 15000
$ws.Cells.Item(15, 13).Value = 14500
$ws.Cells.Item(15, 14).Value = "$/bandeja 18 kilos"
$ws.Cells.Item(15, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(15, 16).Value = 806
$ws.Cells.Item(15, 17).Value = 18
$ws.Cells.Item(15, 18).Value = "Hortaliza"

# New row 16 - "Segunda" quality, this week's data
$ws.Cells.Item(16, 1).Value = 11
$ws.Cells.Item(16, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(16, 3).Value = "Bíobío"
$ws.Cells.Item(16, 4).Value = 44742
$ws.Cells.Item(16, 5).Value = 8
$ws.Cells.Item(16, 6).Value = 100112043
$ws.Cells.Item(16, 7).Value = "Pepino dulce"
$ws.Cells.Item(16, 8).Value = "Cultivar IV Región"
$ws.Cells.Item(16, 9).Value = "Segunda"
$ws.Cells.Item(16, 10).Value = 50
$ws.Cells.Item(16, 11).Value = 12000
$ws.Cells.Item(16, 12).Value = 12000
$ws.Cells.Item(16, 13).Value = 12000
$ws.Cells.Item(16, 14).Value = "$/bandeja 18 kilos"
$ws.Cells.Item(16, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(16, 16).Value = 667
$ws.Cells.Item(16, 17).Value = 18
$ws.Cells.Item(16, 18).Value = "Hortaliza"
